$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table currently holds years 2000-2020 in rows 2-22.
# We need the table to hold years 2010-2022 in rows 2-14:
#   - drop the oldest 10 years (2000-2009, rows 2-11), shifting the rest up
#   - append two new years (2021, 2022) at the bottom

# Remove rows for 2000-2009 (old rows 2 through 11); remaining rows shift up.
$ws.Range("A2:F11").EntireRow.Delete()

# After the delete, rows 2-12 hold 2010-2020. Append 2021 data in row 13.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 22957.9281
$ws.Range("C13").Value = 1500.0319
$ws.Range("D13").Value = 1459.5954
$ws.Range("E13").Value = 3312.5261
$ws.Range("F13").Value = 258.7534

# Append 2022 data in row 14 (only total participants and benefit recipients known).
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 23807
$ws.Range("F14").Value = 297

# Match the formatting used by the rest of column A (bold, centered, bordered).
$ws.Range("A2").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
